# Tue, Jul 07, 2020 11:05:39 AM
#
# The deck currently uses the "Integral" (Red Violet) theme on its single
# slide master (ppt/theme/theme2.xml); the Notes Master carries the
# default "Office Theme" colours (ppt/theme/theme1.xml), unused anywhere
# else. The edit swaps the two themes around, so the slides themselves
# go back to the plain default "Office" colour scheme.
#
# PowerPoint's object model doesn't give a "swap theme parts" verb, so we
# reproduce the effect by rewriting every slot of the active theme's
# colour scheme (reachable through Slide.ThemeColorScheme, which is
# backed by the one theme part every slide/layout/master shares) to the
# stock Office palette values.

function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$theme = $s.ThemeColorScheme

$theme.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$theme.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$theme.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$theme.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$theme.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$theme.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$theme.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$theme.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$theme.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$theme.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$theme.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$theme.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
